# Auto-generated Excel COM-interop script
# Applies scheduled market-price data refresh values to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, LTW, WVR) per the commit's canonical OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 11112217
$ws.Range("I92").Value = 13334184
$ws.Range("K92").Value = 13334184
$ws.Range("M92").Value = -13332936
# Row 131
$ws.Range("H131").Value = 2054.04
$ws.Range("I131").Value = 1030.3846
$ws.Range("J131").Value = 3163
$ws.Range("K131").Value = 3091.1538
$ws.Range("L131").Value = 9489
$ws.Range("M131").Value = 1948.8462
$ws.Range("N131").Value = -19569
# Row 132
$ws.Range("H132").Value = 1108.6809
$ws.Range("I132").Value = 713.6222
$ws.Range("J132").Value = 9997.5
$ws.Range("K132").Value = 2140.8666
$ws.Range("L132").Value = 29992.5
$ws.Range("M132").Value = 389.1333999999997
$ws.Range("N132").Value = -35052.5
# Row 137
$ws.Range("H137").Value = 2129.25
$ws.Range("I137").Value = 1930.6666
$ws.Range("J137").Value = 2725
$ws.Range("K137").Value = 5791.9998
$ws.Range("L137").Value = 8175
$ws.Range("M137").Value = -3241.9998
$ws.Range("N137").Value = -13275
# Row 138
$ws.Range("H138").Value = 2057.3425
$ws.Range("I138").Value = 1388.1951
$ws.Range("J138").Value = 2914.6875
$ws.Range("K138").Value = 4164.5853
$ws.Range("L138").Value = 8744.0625
$ws.Range("M138").Value = 975.4147000000003
$ws.Range("N138").Value = -19024.0625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 648602.1
$ws.Range("I32").Value = 704451.1
$ws.Range("K32").Value = 704451.1
$ws.Range("M32").Value = -704164.1
# Row 61
$ws.Range("H61").Value = 1842.7076
$ws.Range("I61").Value = 1599.12
$ws.Range("J61").Value = 2654.6667
$ws.Range("K61").Value = 1599.12
$ws.Range("L61").Value = 2654.6667
$ws.Range("M61").Value = -1387.12
$ws.Range("N61").Value = -3078.6667
# Row 74
$ws.Range("H74").Value = 1598.8431
$ws.Range("I74").Value = 1303.2222
$ws.Range("K74").Value = 1303.2222
$ws.Range("M74").Value = -429.2221999999999
# Row 77
$ws.Range("H77").Value = 1598.8431
$ws.Range("I77").Value = 1303.2222
$ws.Range("K77").Value = 6516.111
$ws.Range("M77").Value = -2148.111
# Row 136
$ws.Range("H136").Value = 1842.7076
$ws.Range("I136").Value = 1599.12
$ws.Range("J136").Value = 2654.6667
$ws.Range("K136").Value = 4797.36
$ws.Range("L136").Value = 7964.000100000001
$ws.Range("M136").Value = -2247.36
$ws.Range("N136").Value = -13064.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3221.9614
$ws.Range("I134").Value = 2355.9473
$ws.Range("J134").Value = 5572.5713
$ws.Range("K134").Value = 7067.841899999999
$ws.Range("L134").Value = 16717.7139
$ws.Range("M134").Value = -4532.841899999999
$ws.Range("N134").Value = -21787.7139

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 106.05882
$ws.Range("I7").Value = 96.30768999999999
$ws.Range("J7").Value = 137.75
$ws.Range("K7").Value = 96.30768999999999
$ws.Range("L7").Value = 137.75
$ws.Range("M7").Value = 16.69231000000001
$ws.Range("N7").Value = -363.75
# Row 31
$ws.Range("H31").Value = 5266
$ws.Range("I31").Value = 1052.7097
$ws.Range("J31").Value = 9223.939
$ws.Range("K31").Value = 1052.7097
$ws.Range("L31").Value = 9223.939
$ws.Range("M31").Value = -757.7097000000001
$ws.Range("N31").Value = -9813.939
# Row 34
$ws.Range("H34").Value = 5266
$ws.Range("I34").Value = 1052.7097
$ws.Range("J34").Value = 9223.939
$ws.Range("K34").Value = 1052.7097
$ws.Range("L34").Value = 9223.939
$ws.Range("M34").Value = -850.7097000000001
$ws.Range("N34").Value = -9627.939
# Row 99
$ws.Range("H99").Value = 1923.8096
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1923.8096
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1923.8096
$ws.Range("N99").Value = -4919.809600000001
$ws.Range("M99").ClearContents()
# Row 105
$ws.Range("H105").Value = 2460.5
$ws.Range("I105").Value = 910
$ws.Range("J105").Value = 4011
$ws.Range("K105").Value = 910
$ws.Range("L105").Value = 4011
$ws.Range("M105").Value = 837
$ws.Range("N105").Value = -7505
# Row 126
$ws.Range("H126").Value = 1923.8096
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1923.8096
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5771.4288
$ws.Range("N126").Value = -10711.4288
$ws.Range("M126").ClearContents()
# Row 132
$ws.Range("H132").Value = 3969785.5
$ws.Range("I132").Value = 1226.2
$ws.Range("J132").Value = 23812582
$ws.Range("K132").Value = 3678.6
$ws.Range("L132").Value = 71437746
$ws.Range("M132").Value = -1148.6
$ws.Range("N132").Value = -71442806
# Row 134
$ws.Range("H134").Value = 5977.56
$ws.Range("I134").Value = 6010.6
$ws.Range("J134").Value = 5845.4
$ws.Range("K134").Value = 18031.8
$ws.Range("L134").Value = 17536.2
$ws.Range("M134").Value = -15496.8
$ws.Range("N134").Value = -22606.2

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 5926.35
$ws.Range("I122").Value = 592.7778
$ws.Range("J122").Value = 10290.182
$ws.Range("K122").Value = 5335.000199999999
$ws.Range("L122").Value = 92611.63800000001
$ws.Range("M122").Value = -2885.000199999999
$ws.Range("N122").Value = -97511.63800000001
# Row 131
$ws.Range("H131").Value = 1126.5385
$ws.Range("J131").Value = 1178.75
$ws.Range("L131").Value = 3536.25
$ws.Range("N131").Value = -13616.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 12756.056
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 13406.412
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 13406.412
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -13996.412
# Row 27
$ws.Range("H27").Value = 12756.056
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 13406.412
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 13406.412
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -13620.412
# Row 122
$ws.Range("H122").Value = 3146.0264
$ws.Range("I122").Value = 1752.2
$ws.Range("J122").Value = 4694.722
$ws.Range("K122").Value = 5256.6
$ws.Range("L122").Value = 14084.166
$ws.Range("M122").Value = -2806.6
$ws.Range("N122").Value = -18984.166

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2875387.8
$ws.Range("I132").Value = 1778.4736
$ws.Range("J132").Value = 8335245
$ws.Range("K132").Value = 5335.4208
$ws.Range("L132").Value = 25005735
$ws.Range("M132").Value = -2805.4208
$ws.Range("N132").Value = -25010795

